$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting "Ano" and "Índice" one column to the right
$ws.Range("B1").EntireColumn.Insert()

# New header for the inserted column (matches the bold/centered style of the other header cells)
$ws.Range("B1").Value = "Variável"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108

# Fill the new column with the constant variable label for every data row (2-31)
$ws.Range("B2:B31").Value = "Índice do emprego formal: 2012=100"
